# Add SVR parameter loading from pred_par structure and Excel files
# New columns K, L, M: svr_kernel_scale, svr_epsilon, svr_box_constraint

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers (row 1) ---
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# --- New parameter values (row 2) ---
$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 0.05
$ws.Range("M2").Value = 20

# --- Drop the now-unused leftover formatting on the blank helper row (row 13) ---
$ws.Range("A13").Style = "Normal"
$ws.Range("A13").ClearContents()

# --- Style table clean-up: cells that only carried the (duplicate/default)
#     "Normal" style lose their explicit style so they fall back to the
#     sheet default, matching the compacted style table written on save ---
$normalRanges = @("B1","C1","H1","H2","A5","A6","A7","A8","A11","D14:G14","D15:G15","I15","D16:G16","I16")
foreach ($addr in $normalRanges) {
    $ws.Range($addr).Style = "Normal"
}

# row 7 also loses its row-level custom formatting
$ws.Rows("7:7").ClearFormats()

# --- Selection moved from J8 to I8 ---
$ws.Range("I8").Select() | Out-Null
